$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Stories")

# Header row for the new PERT columns
$ws.Range("J1").Value = "Story ID (PERT)"
$ws.Range("K1").Value = "Predecessor"
$ws.Range("J1:K1").Font.Bold = $true

# Story IDs (column J) and Predecessors (column K) for Sprint 1 stories (rows 2-13)
$ws.Range("J2").Value = "S1-01"
$ws.Range("K2").Value = "-"

$ws.Range("J3").Value = "S1-02"
$ws.Range("K3").Value = "S1-12"

$ws.Range("J4").Value = "S1-03"
$ws.Range("K4").Value = "S1-12"

$ws.Range("J5").Value = "S1-04"
$ws.Range("K5").Value = "S1-12"

$ws.Range("J6").Value = "S1-05"
$ws.Range("K6").Value = "S1-03"

$ws.Range("J7").Value = "S1-06"
$ws.Range("K7").Value = "S1-05"

$ws.Range("J8").Value = "S1-07"
$ws.Range("K8").Value = "S1-09, 02, 06, 08"

$ws.Range("J9").Value = "S1-08"
$ws.Range("K9").Value = "S1-05"

$ws.Range("J10").Value = "S1-09"
$ws.Range("K10").Value = "S1-04"

$ws.Range("J11").Value = "S1-10"
$ws.Range("K11").Value = "S1-03"

$ws.Range("J12").Value = "S1-11"
$ws.Range("K12").Value = "S1-07, 10"

$ws.Range("J13").Value = "S1-12"
$ws.Range("K13").Value = "S1-01"

# Column widths for the new columns
$ws.Columns.Item(10).ColumnWidth = 14.28515625
$ws.Columns.Item(11).ColumnWidth = 16

# Row 10's Expected (Hours) formula becomes a non-shared explicit formula
$ws.Range("I10").Formula = "=((F10+(4*G10)+H10)/6)"

# Selection moves to J2
$ws.Range("J2").Select()

Write-Host "done"
